$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue "D2" "68.703.61"
Set-TextValue "E2" "  -0.47%  "
Set-TextValue "D3" "3.751.60"
Set-TextValue "E3" "  -1.75%  "
Set-TextValue "E4" "  +0.12%  "
Set-TextValue "D5" "627.23"
Set-TextValue "E5" "  -0.42%  "
Set-TextValue "D6" "164.73"
Set-TextValue "E6" "  -0.31%  "
Set-TextValue "D7" "3.749.78"
Set-TextValue "E7" "  -1.68%  "
Set-TextValue "E8" "  -0.08%  "
Set-TextValue "E9" "  +0.20%  "
Set-TextValue "D10" "0.157"
Set-TextValue "E10" "  -2.54%  "
Set-TextValue "E11" "  -0.06%  "
Set-TextValue "E12" "  +4.13%  "
Set-TextValue "D13" "0.0000237"
Set-TextValue "E13" "  -5.02%  "
Set-TextValue "D14" "34.68"
Set-TextValue "E14" "  -3.73%  "
Set-TextValue "D15" "4.387.65"
Set-TextValue "E15" "  -1.59%  "
Set-TextValue "D16" "3.752.07"
Set-TextValue "E16" "  +0.17%  "
Set-TextValue "D17" "68.699.46"
Set-TextValue "E17" "  -0.53%  "
Set-TextValue "D18" "17.53"
Set-TextValue "E18" "  -2.55%  "
Set-TextValue "E19" "  -0.43%  "
Set-TextValue "D20" "6.96"
Set-TextValue "E20" "  -2.28%  "
Set-TextValue "D21" "468.26"
Set-TextValue "E21" "  +0.63%  "
Set-TextValue "D22" "9.43"
Set-TextValue "E22" "  -2.75%  "
Set-TextValue "D23" "0.699"
Set-TextValue "E23" "  -1.24%  "
Set-TextValue "D24" "81.47"
Set-TextValue "D25" "0.0000141"
Set-TextValue "E25" "  -6.95%  "
Set-TextValue "E27" "  -2.84%  "
Set-TextValue "D28" "10.01"
Set-TextValue "E28" "  -0.14%  "
Set-TextValue "E29" "  -0.09%  "
Set-TextValue "D30" "3.902.07"
Set-TextValue "E30" "  -1.65%  "
Set-TextValue "D31" "2.25"
Set-TextValue "E31" "  +1.42%  "
Set-TextValue "E32" "  -1.81%  "
Set-TextValue "D33" "7.08"
Set-TextValue "E33" "  -2.78%  "
Set-TextValue "E34" "  +19.18%  "
Set-TextValue "E36" "  +0.10%  "
Set-TextValue "D37" "3.706.53"
Set-TextValue "E37" "  -1.50%  "
Set-TextValue "D38" "8.83"
Set-TextValue "E38" "  -2.80%  "
Set-TextValue "E39" "  -1.28%  "
Set-TextValue "E40" "  -5.42%  "
Set-TextValue "D41" "5.73"
Set-TextValue "E42" "  +0.03%  "
Set-TextValue "D43" "0.953"
Set-TextValue "E43" "  -2.82%  "
Set-TextValue "D45" "44.33"
Set-TextValue "E45" "  +4.46%  "
Set-TextValue "D46" "156.27"
Set-TextValue "E46" "  -0.67%  "
Set-TextValue "D48" "47.11"
Set-TextValue "E48" "  +0.65%  "
Set-TextValue "E49" "  -3.14%  "
Set-TextValue "D50" "0.291"
Set-TextValue "E50" "  -2.73%  "
Set-TextValue "E51" "  -1.59%  "
